$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.261.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.896.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.692"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.39%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.46"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.348"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0722"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0986"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.171.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.708"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.902.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.259.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0820"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "240.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.86%  "
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.130"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.24%  "
$ws.Range("E31").Value = "  +20.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.49%  "
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.37%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +13.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.911"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0659"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.23%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.96%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "93.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.50%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0207"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.349.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.07%  "
